# Weekly update: insert a new price record at row 20 for "Vega Modelo de
# Temuco - Arveja Verde". Existing rows 20-50 shift down to 21-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20, pushing rows 20-50 down to 21-51.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new market record.
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 44498
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 100112022
$ws.Range("G20").Value = "Arveja Verde"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("M20").Value = 20000
$ws.Range("N20").Value = '$/malla 25 kilos'
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 800
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"
